$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2749.5
$ws.Range("I62").Value = 2749.5
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 2749.5
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -2125.5
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 2749.5
$ws.Range("I65").Value = 2749.5
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 13747.5
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -10627.5
$ws.Range("N65").ClearContents()

$ws.Range("H129").Value = 934.2714
$ws.Range("I129").Value = 366.33334
$ws.Range("J129").Value = 959.7015
$ws.Range("K129").Value = 1099.00002
$ws.Range("L129").Value = 2879.1045
$ws.Range("M129").Value = 3900.99998
$ws.Range("N129").Value = -12879.1045

$ws.Range("H134").Value = 59070.715
$ws.Range("J134").Value = 59070.715
$ws.Range("L134").Value = 59070.715
$ws.Range("N134").Value = -69210.715

$ws.Range("H136").Value = 58317.273
$ws.Range("J136").Value = 58317.273
$ws.Range("L136").Value = 58317.273
$ws.Range("N136").Value = -68517.273

$ws.Range("H138").Value = 2916.202
$ws.Range("I138").Value = 1701.2273
$ws.Range("J138").Value = 3263.3376
$ws.Range("K138").Value = 5103.6819
$ws.Range("L138").Value = 9790.0128
$ws.Range("M138").Value = 36.31810000000041
$ws.Range("N138").Value = -20070.0128

$ws.Range("H139").Value = 38931.305
$ws.Range("J139").Value = 38931.305
$ws.Range("L139").Value = 38931.305
$ws.Range("N139").Value = -49211.305

$ws.Range("H140").Value = 49234.484
$ws.Range("J140").Value = 49234.484
$ws.Range("L140").Value = 49234.484
$ws.Range("N140").Value = -59594.484

$ws.Range("H141").Value = 5838.423
$ws.Range("I141").Value = 6077.227
$ws.Range("J141").Value = 4525
$ws.Range("K141").Value = 18231.681
$ws.Range("L141").Value = 13575
$ws.Range("M141").Value = -13051.681
$ws.Range("N141").Value = -23935

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5995.83
$ws.Range("I32").Value = 3503.7974
$ws.Range("J32").Value = 13088.538
$ws.Range("K32").Value = 3503.7974
$ws.Range("L32").Value = 13088.538
$ws.Range("M32").Value = -3216.7974
$ws.Range("N32").Value = -13662.538

$ws.Range("H132").Value = 2219.682
$ws.Range("I132").Value = 1096.5
$ws.Range("J132").Value = 5214.8335
$ws.Range("K132").Value = 3289.5
$ws.Range("L132").Value = 15644.5005
$ws.Range("M132").Value = -759.5
$ws.Range("N132").Value = -20704.5005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 259.96155
$ws.Range("I80").Value = 115.545456
$ws.Range("J80").Value = 365.86667
$ws.Range("K80").Value = 115.545456
$ws.Range("L80").Value = 365.86667
$ws.Range("M80").Value = 882.4545439999999
$ws.Range("N80").Value = -2361.86667

$ws.Range("H83").Value = 259.96155
$ws.Range("I83").Value = 115.545456
$ws.Range("J83").Value = 365.86667
$ws.Range("K83").Value = 577.7272800000001
$ws.Range("L83").Value = 1829.33335
$ws.Range("M83").Value = 4414.27272
$ws.Range("N83").Value = -11813.33335

$ws.Range("H100").Value = 26000
$ws.Range("J100").Value = 26000
$ws.Range("L100").Value = 26000
$ws.Range("N100").Value = -28164

$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

$ws.Range("H109").Value = 38000
$ws.Range("J109").Value = 38000
$ws.Range("L109").Value = 38000
$ws.Range("N109").Value = -40774

$ws.Range("H122").Value = 43112
$ws.Range("J122").Value = 43112
$ws.Range("L122").Value = 43112
$ws.Range("N122").Value = -52912

$ws.Range("H126").Value = 43780
$ws.Range("J126").Value = 43780
$ws.Range("L126").Value = 43780
$ws.Range("N126").Value = -53660

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 730
$ws.Range("I107").Value = 746.6667
$ws.Range("J107").Value = 700
$ws.Range("K107").Value = 746.6667
$ws.Range("L107").Value = 700
$ws.Range("M107").Value = 1173.3333
$ws.Range("N107").Value = -4540

$ws.Range("H132").Value = 2111.7737
$ws.Range("I132").Value = 1698.7954
$ws.Range("J132").Value = 4130.778
$ws.Range("K132").Value = 5096.3862
$ws.Range("L132").Value = 12392.334
$ws.Range("M132").Value = -2566.3862
$ws.Range("N132").Value = -17452.334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 11633664
$ws.Range("I131").Value = 55579388
$ws.Range("J131").Value = 971.58826
$ws.Range("K131").Value = 166738164
$ws.Range("L131").Value = 2914.76478
$ws.Range("M131").Value = -166733124
$ws.Range("N131").Value = -12994.76478

$ws.Range("H132").Value = 5374.1875
$ws.Range("I132").Value = 1411.7142
$ws.Range("K132").Value = 12705.4278
$ws.Range("M132").Value = -10175.4278

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H108").Value = 50000
$ws.Range("J108").Value = 50000
$ws.Range("L108").Value = 50000
$ws.Range("N108").Value = -57680

$ws.Range("H132").Value = 2661.725
$ws.Range("I132").Value = 1148.7273
$ws.Range("J132").Value = 3235.6206
$ws.Range("K132").Value = 3446.1819
$ws.Range("L132").Value = 9706.861800000001
$ws.Range("M132").Value = -916.1819
$ws.Range("N132").Value = -14766.8618

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5260.722
$ws.Range("I7").Value = 3351.1
$ws.Range("K7").Value = 3351.1
$ws.Range("M7").Value = -3239.1

$ws.Range("H16").Value = 482.83334
$ws.Range("I16").Value = 500.5
$ws.Range("J16").Value = 235.5
$ws.Range("K16").Value = 500.5
$ws.Range("L16").Value = 235.5
$ws.Range("M16").Value = -330.5
$ws.Range("N16").Value = -575.5

$ws.Range("H40").Value = 4912.0566
$ws.Range("I40").Value = 4592.07
$ws.Range("K40").Value = 4592.07
$ws.Range("M40").Value = -4456.07

$ws.Range("H56").Value = 16721.334
$ws.Range("I56").Value = 10051
$ws.Range("J56").Value = 20056.5
$ws.Range("K56").Value = 10051
$ws.Range("L56").Value = 20056.5
$ws.Range("M56").Value = -9360
$ws.Range("N56").Value = -21438.5

$ws.Range("H126").Value = 5260.722
$ws.Range("I126").Value = 3351.1
$ws.Range("K126").Value = 10053.3
$ws.Range("M126").Value = -7583.299999999999

$ws.Range("H139").Value = 47136.668
$ws.Range("J139").Value = 50705
$ws.Range("L139").Value = 50705
$ws.Range("N139").Value = -60985

$ws.Range("H140").Value = 75000
$ws.Range("J140").Value = 75000
$ws.Range("L140").Value = 75000
$ws.Range("N140").Value = -85360

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 9264047
$ws.Range("I132").Value = 6709.9414
$ws.Range("K132").Value = 20129.8242
$ws.Range("M132").Value = -17599.8242

$ws.Range("H136").Value = 3330.2258
$ws.Range("I136").Value = 937.7619
$ws.Range("J136").Value = 8354.4
$ws.Range("K136").Value = 2813.2857
$ws.Range("L136").Value = 25063.2
$ws.Range("M136").Value = -263.2856999999999
$ws.Range("N136").Value = -30163.2
